$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 33.64214033333334
$ws.Range("H2").Value = 100.926421
$ws.Range("I2").Value = 0.106995191696894
$ws.Range("J2").Value = 0.106995191696894
$ws.Range("M2").Value = 0.01339666666666667
$ws.Range("N2").Value = 0.04019
$ws.Range("O2").Value = 0.08393217762128816
$ws.Range("P2").Value = 0.08393217762128814
$ws.Range("Q2").Value = 0.4506925399988889
$ws.Range("R2").Value = 4.05623285999
$ws.Range("S2").Value = 0.008980339434127483
$ws.Range("T2").Value = 0.008980339434127479
$ws.Range("G3").Value = 33.64214033333334
$ws.Range("H3").Value = 100.926421
$ws.Range("I3").Value = 0.106995191696894
$ws.Range("J3").Value = 0.106995191696894
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.1260863333333333
$ws.Range("N3").Value = 0.378259
$ws.Range("O3").Value = 0.7899502755623498
$ws.Range("P3").Value = 0.7899502755623498
$ws.Range("Q3").Value = 4.241814120115444
$ws.Range("R3").Value = 38.17632708103901
$ws.Range("S3").Value = 0.08452088116480785
$ws.Range("T3").Value = 0.08452088116480784
$ws.Range("G4").Value = 33.64214033333334
$ws.Range("H4").Value = 100.926421
$ws.Range("I4").Value = 0.106995191696894
$ws.Range("J4").Value = 0.106995191696894
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02013
$ws.Range("N4").Value = 0.06039
$ws.Range("O4").Value = 0.1261175468163621
$ws.Range("P4").Value = 0.1261175468163621
$ws.Range("Q4").Value = 0.6772162849100001
$ws.Range("R4").Value = 6.09494656419
$ws.Range("S4").Value = 0.01349397109795866
$ws.Range("T4").Value = 0.01349397109795866
$ws.Range("H5").Value = 727.0751789999999
$ws.Range("I5").Value = 0.7707946777896593
$ws.Range("J5").Value = 0.7707946777896592
$ws.Range("M5").Value = 0.01339666666666667
$ws.Range("N5").Value = 0.04019
$ws.Range("O5").Value = 0.08393217762128816
$ws.Range("P5").Value = 0.08393217762128814
$ws.Range("Q5").Value = 3.246794604889999
$ws.Range("R5").Value = 29.22115144400999
$ws.Range("S5").Value = 0.06469447580578526
$ws.Range("T5").Value = 0.06469447580578523
$ws.Range("H6").Value = 727.0751789999999
$ws.Range("I6").Value = 0.7707946777896593
$ws.Range("J6").Value = 0.7707946777896592
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.1260863333333333
$ws.Range("N6").Value = 0.378259
$ws.Range("O6").Value = 0.7899502755623498
$ws.Range("P6").Value = 0.7899502755623498
$ws.Range("Q6").Value = 30.558081125929
$ws.Range("R6").Value = 275.022730133361
$ws.Range("S6").Value = 0.6088894681219339
$ws.Range("T6").Value = 0.6088894681219339
$ws.Range("H7").Value = 727.0751789999999
$ws.Range("I7").Value = 0.7707946777896593
$ws.Range("J7").Value = 0.7707946777896592
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.02013
$ws.Range("N7").Value = 0.06039
$ws.Range("O7").Value = 0.1261175468163621
$ws.Range("P7").Value = 0.1261175468163621
$ws.Range("Q7").Value = 4.878674451089999
$ws.Range("R7").Value = 43.90807005980999
$ws.Range("S7").Value = 0.09721073386194005
$ws.Range("T7").Value = 0.09721073386194004
$ws.Range("G8").Value = 9.788214000000002
$ws.Range("H8").Value = 29.364642
$ws.Range("I8").Value = 0.0311303568359039
$ws.Range("J8").Value = 0.03113035683590389
$ws.Range("M8").Value = 0.01339666666666667
$ws.Range("N8").Value = 0.04019
$ws.Range("O8").Value = 0.08393217762128816
$ws.Range("P8").Value = 0.08393217762128814
$ws.Range("Q8").Value = 0.13112944022
$ws.Range("R8").Value = 1.18016496198
$ws.Range("S8").Value = 0.002612838639365168
$ws.Range("T8").Value = 0.002612838639365167
$ws.Range("G9").Value = 9.788214000000002
$ws.Range("H9").Value = 29.364642
$ws.Range("I9").Value = 0.0311303568359039
$ws.Range("J9").Value = 0.03113035683590389
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.1260863333333333
$ws.Range("N9").Value = 0.378259
$ws.Range("O9").Value = 0.7899502755623498
$ws.Range("P9").Value = 0.7899502755623498
$ws.Range("Q9").Value = 1.234160013142
$ws.Range("R9").Value = 11.107440118278
$ws.Range("S9").Value = 0.02459143396087657
$ws.Range("T9").Value = 0.02459143396087656
$ws.Range("G10").Value = 9.788214000000002
$ws.Range("H10").Value = 29.364642
$ws.Range("I10").Value = 0.0311303568359039
$ws.Range("J10").Value = 0.03113035683590389
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.02013
$ws.Range("N10").Value = 0.06039
$ws.Range("O10").Value = 0.1261175468163621
$ws.Range("P10").Value = 0.1261175468163621
$ws.Range("Q10").Value = 0.19703674782
$ws.Range("R10").Value = 1.77333073038
$ws.Range("S10").Value = 0.003926084235662167
$ws.Range("T10").Value = 0.003926084235662166
$ws.Range("G11").Value = 12.12016933333333
$ws.Range("H11").Value = 36.360508
$ws.Range("I11").Value = 0.03854688876420623
$ws.Range("J11").Value = 0.03854688876420622
$ws.Range("M11").Value = 0.01339666666666667
$ws.Range("N11").Value = 0.04019
$ws.Range("O11").Value = 0.08393217762128816
$ws.Range("P11").Value = 0.08393217762128814
$ws.Range("Q11").Value = 0.1623698685022222
$ws.Range("R11").Value = 1.46132881652
$ws.Range("S11").Value = 0.003235324314505394
$ws.Range("T11").Value = 0.003235324314505393
$ws.Range("G12").Value = 12.12016933333333
$ws.Range("H12").Value = 36.360508
$ws.Range("I12").Value = 0.03854688876420623
$ws.Range("J12").Value = 0.03854688876420622
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.1260863333333333
$ws.Range("N12").Value = 0.378259
$ws.Range("O12").Value = 0.7899502755623498
$ws.Range("P12").Value = 0.7899502755623498
$ws.Range("Q12").Value = 1.528187710619111
$ws.Range("R12").Value = 13.753689395572
$ws.Range("S12").Value = 0.03045012540135596
$ws.Range("T12").Value = 0.03045012540135595
$ws.Range("G13").Value = 12.12016933333333
$ws.Range("H13").Value = 36.360508
$ws.Range("I13").Value = 0.03854688876420623
$ws.Range("J13").Value = 0.03854688876420622
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.02013
$ws.Range("N13").Value = 0.06039
$ws.Range("O13").Value = 0.1261175468163621
$ws.Range("P13").Value = 0.1261175468163621
$ws.Range("Q13").Value = 0.24397900868
$ws.Range("R13").Value = 2.19581107812
$ws.Range("S13").Value = 0.004861439048344881
$ws.Range("T13").Value = 0.004861439048344879
$ws.Range("G14").Value = 16.51773933333333
$ws.Range("H14").Value = 49.553218
$ws.Range("I14").Value = 0.0525328849133368
$ws.Range("J14").Value = 0.05253288491333678
$ws.Range("M14").Value = 0.01339666666666667
$ws.Range("N14").Value = 0.04019
$ws.Range("O14").Value = 0.08393217762128816
$ws.Range("P14").Value = 0.08393217762128814
$ws.Range("Q14").Value = 0.2212826479355556
$ws.Range("R14").Value = 1.99154383142
$ws.Range("S14").Value = 0.004409199427504873
$ws.Range("T14").Value = 0.004409199427504871
$ws.Range("G15").Value = 16.51773933333333
$ws.Range("H15").Value = 49.553218
$ws.Range("I15").Value = 0.0525328849133368
$ws.Range("J15").Value = 0.05253288491333678
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 0.3333333333333333
$ws.Range("M15").Value = 0.1260863333333333
$ws.Range("N15").Value = 0.378259
$ws.Range("O15").Value = 0.7899502755623498
$ws.Range("P15").Value = 0.7899502755623498
$ws.Range("Q15").Value = 2.082661187495778
$ws.Range("R15").Value = 18.743950687462
$ws.Range("S15").Value = 0.04149836691337561
$ws.Range("T15").Value = 0.0414983669133756
$ws.Range("G16").Value = 16.51773933333333
$ws.Range("H16").Value = 49.553218
$ws.Range("I16").Value = 0.0525328849133368
$ws.Range("J16").Value = 0.05253288491333678
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.02013
$ws.Range("N16").Value = 0.06039
$ws.Range("O16").Value = 0.1261175468163621
$ws.Range("P16").Value = 0.1261175468163621
$ws.Range("Q16").Value = 0.33250209278
$ws.Range("R16").Value = 2.99251883502
$ws.Range("S16").Value = 0.006625318572456314
$ws.Range("T16").Value = 0.006625318572456312
